# Dhali's Menu update — revised Thali / Sandwich / Pizza descriptions and prices.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Updated Thali descriptions (rows 13-14)
$ws1.Range("A13").Value = "Veg Thail(Chole, Mix Veg, 2 Roti, Rice, Raita, salad, Aachar, One Sweet)"
$ws1.Range("A14").Value = "Veg Special Thail(Shai Paneer, Mix Veg, 2 Roti, Rice, Raita, salad,Aachar, One Sweet)"

# Updated Schezwan Sandwich description (row 11)
$ws1.Range("A11").Value = "Schezwan Grilled Sandwich – Indo-Chinese fusion with Schezwan sauce, Veggies Patty, and cheese."

# Updated Club Sandwich description (row 12) and its Half price
$ws1.Range("A12").Value = "Club Sandwich (Indian Style) –Paneer Patty, and cheese, Veggies, green chutney"
$ws1.Range("C12").Value = 130

# Updated Pizza names/sizes (rows 5-7)
$ws1.Range("A5").Value = "Pizza Margarita 10' "
$ws1.Range("A6").Value = "Pizza Onion and Capsicum 10' "
$ws1.Range("A7").Value = "Pizza Panner, Veggie ( Onion and Capsicum and corn) 10'"

# Update the saved selection to match the author's last-clicked cell
$ws1.Activate()
$ws1.Range("A4").Select()
